$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1392.4667
$ws.Range("I58").Value = 218.14285
$ws.Range("J58").Value = 2420
$ws.Range("K58").Value = 654.4285500000001
$ws.Range("L58").Value = 7260
$ws.Range("M58").Value = -504.4285500000001
$ws.Range("N58").Value = -7560

$ws.Range("H93").Value = 38000
$ws.Range("J93").Value = 38000
$ws.Range("L93").Value = 38000
$ws.Range("N93").Value = -42992

$ws.Range("H129").Value = 1039.3062
$ws.Range("J129").Value = 1145.7317
$ws.Range("L129").Value = 3437.1951
$ws.Range("N129").Value = -13437.1951

$ws.Range("H140").Value = 56272.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 56272.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 56272.5
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -66632.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9229.218000000001
$ws.Range("I32").Value = 6906.7095
$ws.Range("K32").Value = 6906.7095
$ws.Range("M32").Value = -6619.7095

$ws.Range("H45").Value = 7976
$ws.Range("I45").Value = 9091.538
$ws.Range("J45").Value = 725
$ws.Range("K45").Value = 9091.538
$ws.Range("L45").Value = 725
$ws.Range("M45").Value = -8714.538
$ws.Range("N45").Value = -1479

$ws.Range("H61").Value = 2790.6667
$ws.Range("I61").Value = 1951.5
$ws.Range("J61").Value = 4469
$ws.Range("K61").Value = 1951.5
$ws.Range("L61").Value = 4469
$ws.Range("M61").Value = -1739.5
$ws.Range("N61").Value = -4893

$ws.Range("H74").Value = 1358.4166
$ws.Range("I74").Value = 1022.7879
$ws.Range("J74").Value = 2096.8
$ws.Range("K74").Value = 1022.7879
$ws.Range("L74").Value = 2096.8
$ws.Range("M74").Value = -148.7879
$ws.Range("N74").Value = -3844.8

$ws.Range("H77").Value = 1358.4166
$ws.Range("I77").Value = 1022.7879
$ws.Range("J77").Value = 2096.8
$ws.Range("K77").Value = 5113.9395
$ws.Range("L77").Value = 10484
$ws.Range("M77").Value = -745.9395000000004
$ws.Range("N77").Value = -19220

$ws.Range("H92").Value = 31620
$ws.Range("J92").Value = 31620
$ws.Range("L92").Value = 31620
$ws.Range("N92").Value = -36612

$ws.Range("H123").Value = 39000
$ws.Range("J123").Value = 39000
$ws.Range("L123").Value = 39000
$ws.Range("N123").Value = -48800

$ws.Range("H132").Value = 2705682.2
$ws.Range("I132").Value = 1740.5
$ws.Range("J132").Value = 7697575
$ws.Range("K132").Value = 5221.5
$ws.Range("L132").Value = 23092725
$ws.Range("M132").Value = -2691.5
$ws.Range("N132").Value = -23097785

$ws.Range("H133").Value = 38250
$ws.Range("J133").Value = 38250
$ws.Range("L133").Value = 38250
$ws.Range("N133").Value = -43310

$ws.Range("H136").Value = 2790.6667
$ws.Range("I136").Value = 1951.5
$ws.Range("J136").Value = 4469
$ws.Range("K136").Value = 5854.5
$ws.Range("L136").Value = 13407
$ws.Range("M136").Value = -3304.5
$ws.Range("N136").Value = -18507

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 3268.6667
$ws.Range("I5").Value = 4700.5
$ws.Range("J5").Value = 405
$ws.Range("K5").Value = 4700.5
$ws.Range("L5").Value = 405
$ws.Range("M5").Value = -4587.5
$ws.Range("N5").Value = -631

$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H8").Value = 2000
$ws.Range("I8").Value = 2000
$ws.Range("K8").Value = 2000
$ws.Range("M8").Value = -1860

$ws.Range("H11").Value = 1366.6666
$ws.Range("I11").Value = 550
$ws.Range("K11").Value = 550
$ws.Range("M11").Value = -410

$ws.Range("H12").Value = 1752.5
$ws.Range("I12").Value = 1170
$ws.Range("J12").Value = 3500
$ws.Range("K12").Value = 1170
$ws.Range("L12").Value = 3500
$ws.Range("M12").Value = -1002
$ws.Range("N12").Value = -3836

$ws.Range("H134").Value = 2608.375
$ws.Range("I134").Value = 2393.4
$ws.Range("J134").Value = 2966.6667
$ws.Range("K134").Value = 7180.200000000001
$ws.Range("L134").Value = 8900.000100000001
$ws.Range("M134").Value = -4645.200000000001
$ws.Range("N134").Value = -13970.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5231.115
$ws.Range("I31").Value = 1757.5151
$ws.Range("J31").Value = 11264.211
$ws.Range("K31").Value = 1757.5151
$ws.Range("L31").Value = 11264.211
$ws.Range("M31").Value = -1462.5151
$ws.Range("N31").Value = -11854.211

$ws.Range("H34").Value = 5231.115
$ws.Range("I34").Value = 1757.5151
$ws.Range("J34").Value = 11264.211
$ws.Range("K34").Value = 1757.5151
$ws.Range("L34").Value = 11264.211
$ws.Range("M34").Value = -1555.5151
$ws.Range("N34").Value = -11668.211

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 110231.56
$ws.Range("I4").Value = 250112.6
$ws.Range("J4").Value = 2630.7693
$ws.Range("K4").Value = 750337.8
$ws.Range("L4").Value = 7892.3079
$ws.Range("M4").Value = -750225.8
$ws.Range("N4").Value = -8116.3079

$ws.Range("H21").Value = 1054.4546
$ws.Range("I21").Value = 400
$ws.Range("J21").Value = 1119.9
$ws.Range("K21").Value = 1200
$ws.Range("L21").Value = 3359.7
$ws.Range("M21").Value = -1027
$ws.Range("N21").Value = -3705.7

$ws.Range("H55").Value = 3205.2632
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 3355.5557
$ws.Range("K55").Value = 1500
$ws.Range("L55").Value = 10066.6671
$ws.Range("M55").Value = -1323
$ws.Range("N55").Value = -10420.6671

$ws.Range("H131").Value = 2001044
$ws.Range("J131").Value = 1192.7142
$ws.Range("L131").Value = 3578.1426
$ws.Range("N131").Value = -13658.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 50000
$ws.Range("J4").Value = 50000
$ws.Range("L4").Value = 50000
$ws.Range("N4").Value = -50224

$ws.Range("H40").Value = 12184.615
$ws.Range("J40").Value = 12184.615
$ws.Range("L40").Value = 12184.615
$ws.Range("N40").Value = -12486.615

$ws.Range("H70").Value = 5453.771
$ws.Range("J70").Value = 4824.6
$ws.Range("L70").Value = 4824.6
$ws.Range("N70").Value = -5364.6

$ws.Range("H73").Value = 5453.771
$ws.Range("J73").Value = 4824.6
$ws.Range("L73").Value = 4824.6
$ws.Range("N73").Value = -6696.6

$ws.Range("H132").Value = 3510.1516
$ws.Range("I132").Value = 5284.273
$ws.Range("J132").Value = 2623.0908
$ws.Range("K132").Value = 15852.819
$ws.Range("L132").Value = 7869.2724
$ws.Range("M132").Value = -13322.819
$ws.Range("N132").Value = -12929.2724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 892.1212
$ws.Range("I46").Value = 716.1667
$ws.Range("J46").Value = 992.6667
$ws.Range("K46").Value = 716.1667
$ws.Range("L46").Value = 992.6667
$ws.Range("M46").Value = -528.1667
$ws.Range("N46").Value = -1368.6667

$ws.Range("H136").Value = 4948.881
$ws.Range("I136").Value = 2423.4062
$ws.Range("J136").Value = 13030.4
$ws.Range("K136").Value = 7270.2186
$ws.Range("L136").Value = 39091.2
$ws.Range("M136").Value = -4720.2186
$ws.Range("N136").Value = -44191.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3395.5
$ws.Range("I132").Value = 3101.0833
$ws.Range("J132").Value = 3837.125
$ws.Range("K132").Value = 9303.249899999999
$ws.Range("L132").Value = 11511.375
$ws.Range("M132").Value = -6773.249899999999
$ws.Range("N132").Value = -16571.375

$ws.Range("H136").Value = 2872.9487
$ws.Range("I136").Value = 3612.6316
$ws.Range("J136").Value = 2170.25
$ws.Range("K136").Value = 10837.8948
$ws.Range("L136").Value = 6510.75
$ws.Range("M136").Value = -8287.8948
$ws.Range("N136").Value = -11610.75
